$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (D: Price, E: Volume) are treated as plain text so
# values such as "1.00" or "89.538.33" are not reinterpreted as numbers/dates.
$updates = @{
    "D2" = "89.538.33"
    "E2" = "  +3.57%  "
    "D3" = "3.294.90"
    "E3" = "  -0.84%  "
    "E4" = "  +0.25%  "
    "D5" = "214.92"
    "E5" = "  -1.75%  "
    "D6" = "631.62"
    "E6" = "  -0.68%  "
    "D7" = "0.387"
    "E7" = "  +19.14%  "
    "D8" = "0.689"
    "E8" = "  +15.64%  "
    "D9" = "1.00"
    "E9" = "  +0.11%  "
    "D10" = "3.290.45"
    "E10" = "  -1.00%  "
    "D11" = "0.580"
    "E11" = "  -2.84%  "
    "D12" = "0.188"
    "E12" = "  +12.71%  "
    "E13" = "  -3.65%  "
    "D14" = "34.52"
    "E14" = "  +0.71%  "
    "D15" = "3.887.43"
    "E15" = "  -1.08%  "
    "E16" = "  +0.30%  "
    "D17" = "89.335.98"
    "E17" = "  +4.09%  "
    "D18" = "3.289.04"
    "E18" = "  -0.62%  "
    "B19" = "SuiNetwork"
    "C19" = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
    "D19" = "3.14"
    "E19" = "  -0.97%  "
    "B20" = "Chainlink"
    "C20" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D20" = "14.21"
    "E20" = "  -2.77%  "
    "D21" = "439.65"
    "E21" = "  -0.70%  "
    "D22" = "8.93"
    "E22" = "  -2.29%  "
    "D23" = "5.40"
    "E23" = "  +3.26%  "
    "D24" = "7.38"
    "E24" = "  +0.31%  "
    "D25" = "12.40"
    "E25" = "  +1.03%  "
    "E26" = "  -2.80%  "
    "E27" = "  -0.90%  "
    "D28" = "77.13"
    "E28" = "  -1.54%  "
    "E29" = "  +4.05%  "
    "B30" = "Dai"
    "C30" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D30" = "0.998"
    "E30" = "  -0.20%  "
    "B31" = "Cronos"
    "C31" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D31" = "0.194"
    "E31" = "  +15.77%  "
    "E32" = "  +0.07%  "
    "B33" = "InternetComputer(DFINITY)"
    "C33" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D33" = "8.92"
    "E33" = "  -3.31%  "
    "B34" = "Bittensor"
    "C34" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D34" = "572.12"
    "E34" = "  -6.03%  "
    "E35" = "  -9.51%  "
    "D36" = "7.23"
    "E36" = "  +11.73%  "
    "E37" = "  -3.36%  "
    "E38" = "  -7.39%  "
    "D39" = "22.79"
    "E39" = "  -2.50%  "
    "D40" = "21.86"
    "E40" = "  +2.66%  "
    "D41" = "1.00"
    "E41" = "  +0.35%  "
    "D42" = "3.11"
    "E42" = "  +0.05%  "
    "B43" = "PolygonEcosystemToken"
    "C43" = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
    "D43" = "0.402"
    "E43" = "  -3.44%  "
    "B44" = "Stacks"
    "C44" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D44" = "2.05"
    "E44" = "  -0.54%  "
    "D46" = "155.68"
    "E46" = "  -1.30%  "
    "D47" = "181.99"
    "E47" = "  -3.30%  "
    "D48" = "45.07"
    "E48" = "  -0.66%  "
    "E49" = "  -3.70%  "
    "E50" = "  +22.00%  "
    "B51" = "Stellar"
    "C51" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D51" = "0.127"
    "E51" = "  +15.04%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
